$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.162.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.831.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6601'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.66%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07389'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.67%  '
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07759'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.833.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.990'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6649'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.104'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008403'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.139.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.068.09'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.112'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.589'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.06%  '
$ws.Range("E27").Value = '  -2.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.91'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.27%  '
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.105'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.036'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.189'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05250'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.864'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7390'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("E36").Value = '  +1.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.649'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.295.57'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01788'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.730'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9227'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.961'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9994'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.20'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.971.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5139'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("E48").Value = '  -9.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.745'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05842'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.43%  '
